$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STO LLC Auction Calculations")

# Update existing row 2 amounts to reflect the new auction destination
$ws.Range("H2").Value = 25923.77
$ws.Range("I2").Value = -2656.57

# Add a new row (id 2) for the new auction destination
$ws.Range("A3").Value = 2
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 30000
$ws.Range("B3").HorizontalAlignment = -4131

$ws.Range("C3").Value = "ont"
$ws.Range("C3").HorizontalAlignment = -4131

$ws.Range("D3").Value = 32000
$ws.Range("D3").HorizontalAlignment = -4131

$ws.Range("E3").Value = 23267.2
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = 35030
$ws.Range("F3").Style = "Normal"

$ws.Range("G3").Value = 25470
$ws.Range("G3").Style = "Normal"

$ws.Range("H3").Value = 27138.7
$ws.Range("H3").Style = "Normal"

$ws.Range("I3").Value = -3871.5
$ws.Range("I3").Style = "Normal"

# Restore the recorded selection/active-cell state
$ws.Range("C6").Select() | Out-Null
